$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The worksheet originally ends with a short text note block:
#   A21: "70,2 29,348 "
#   A22..A30: nine more note lines
#
# The edit turns that first note line into a new "doWork02" section header
# (with an extra comment in column B), inserts a brand-new 10-row x 8-column
# results table right under it (mirroring the doWork01 table above), and
# keeps the original note lines (now updated with new numbers, plus two
# extra ones) below that, separated by blank rows - exactly like the
# doWork01 section is laid out earlier in the sheet.
# ---------------------------------------------------------------------------

# 1) Turn row 21 into the new "doWork02" header row, now at row 20. Clear
#    out what is left behind in row 21 so it becomes the blank separator
#    row. (Column B - the long comment - is filled in later, step 4, so
#    that new shared strings get created in the same order the workbook
#    author originally typed them in.)
$ws.Cells.Item(20, 1).Value = "doWork02"
$ws.Cells.Item(21, 1).ClearContents()

# 2) Make room: insert 11 rows at row 22 - ten for the new data table and
#    one blank separator row before the note lines that used to start at
#    row 22 (they will land on row 33 afterwards).
$ws.Rows("22:32").Insert()

# 3) Update the (now shifted) note lines with their new text, and append
#    two brand-new note lines at the end.
$notes = @(
  "doWork03",
  "79,7 41,216 ",
  "  46,1 12,232 ",
  "  40,3  2,094 ",
  "  40,3  1,963 ",
  "  39,8  2,117 ",
  "  39,9  2,769 ",
  "  39,6  1,814 ",
  "  39,7  1,165 ",
  "  40,0  3,097 ",
  "  39,3  1,384 "
)
for ($i = 0; $i -lt $notes.Length; $i++) {
  $row = 33 + $i
  $ws.Cells.Item($row, 1).Value = $notes[$i]
}

# 4) Populate the new doWork02 results table (rows 22-31, columns A-H).
$data = @(
  @(70.2, 29.348, 47.4, 16.309, 57, 38.257, 28.3, 14.875),
  @(39.8, 2.386, 36.3, 1.591, 24.1, 2.47, 21.7, 3.009),
  @(39.3, 1.941, 36.8, 1.23, 25.5, 6.308, 21.7, 2.113),
  @(40.3, 1.841, 36.4, 1.7, 23.3, 1.693, 20.2, 2.043),
  @(39.9, 1.976, 36.1, 1.337, 24, 2.146, 19.6, 0.805),
  @(40.1, 1.612, 37.3, 2.13, 25.2, 3.071, 21.2, 2.038),
  @(39.5, 2.444, 36.8, 1.228, 22.7, 1.435, 20.2, 1.442),
  @(39.3, 1.635, 36.6, 1.319, 21.8, 2.286, 20.3, 1.278),
  @(40, 2.716, 36.4, 1.339, 22.3, 1.896, 21.1, 4.408),
  @(39.2, 1.319, 36.3, 2.546, 22.2, 1.308, 22.3, 6.498)
)
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = 22 + $i
  for ($c = 0; $c -lt 8; $c++) {
    $ws.Cells.Item($row, $c + 1).Value = $data[$i][$c]
  }
}

# 5) Finally fill in the B20 comment, so its shared string is appended
#    last (matching the original authoring order).
$ws.Cells.Item(20, 2).Value = "(byte arrays med 128 eller 256 pladser formentlig bedst.)"

# 6) Match the saved view state: scrolled so row 7 is at the top, with the
#    new doWork02 header cell selected.
$ws.Range("A20").Select()
$excel.ActiveWindow.ScrollRow = 7
